$d = $word.ActiveDocument

# Locate the run-spanning text that needs to change:
#   "Developed front-end/backend " + "tools and " + "technologies for "
# -> "Developed/Integrated" + " front-end/backend " + "tools" +
#    " and technology packages" + " for "
# (the trailing "Mixed Signal Design. " run is left untouched)
$found = $d.Content
$found.Find.Execute("Developed front-end/backend tools and technologies for ", `
                     $true, $false, $false, $false, $false, $true, 1, $false, `
                     "", 0)
$start = $found.Start

# Replace the whole span with the new wording in one go.
$found.Text = "Developed/Integrated front-end/backend tools and technology packages for "

# Re-split that single run back into the six runs implied by the new
# wording by toggling a character-formatting property on/off (a net no-op)
# across each desired sub-range; this forces a run boundary without
# altering any visible formatting.
$parts = @( `
    "Developed/Integrated", `
    " front-end/backend ", `
    "tools", `
    " and technology packages", `
    " for " `
)

$pos = $start
foreach ($part in $parts) {
    $seg = $d.Range($pos, $pos + $part.Length)
    $seg.Font.Bold = 1
    $seg.Font.Bold = 0
    $pos = $pos + $part.Length
}
